# pf2050_data_standard_sample.xlsx
# "Version 0.7; device.device_count added"
#
# Inserts a new row into the "Examples" sheet (row 22) documenting the
# pf2050_data_standard.report.device.device_count field, shifting the
# existing device.substance.* / species.* example rows down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Examples")

# Insert a new blank row at 22 (existing rows 22-34 shift to 23-35),
# inheriting formatting/styles from the row above as Excel normally does.
$ws.Rows("22:22").Insert()

# --- Column A: field name, with "device" bolded (rich text) ---
$ws.Range("A22").Value = "pf2050_data_standard.report.device.device_count"
$fieldName = $ws.Range("A22")
$fieldName.Characters(29, 6).Font.Bold = $true
$fieldName.Characters(35, 13).Font.Bold = $false

# --- Column B: Mandatory/Optional = "O" ---
$ws.Range("B22").Value = "O"

# --- Columns C:J: example counts ---
$ws.Range("C22").Value = 1
$ws.Range("D22").Value = 1
$ws.Range("E22").Value = 1
$ws.Range("F22").Value = 1
$ws.Range("G22").Value = 20
$ws.Range("H22").Value = 20
$ws.Range("I22").Value = 20
$ws.Range("J22").Value = 20

# Row 21 (the row formatting was copied from) has no H/I/J cells, so the
# inserted row doesn't inherit the right-aligned "s=6" style there -
# apply it explicitly to match the rest of the table.
$ws.Range("H22:J22").HorizontalAlignment = -4152

# Move the active selection to C2 (matches the saved workbook state).
$null = $ws.Range("C2").Select()
